$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  17"
$ws.Range("C9").Value = "Report Covering the Week  4/21/2025  Through  4/27/2025"

# --- Crime Complaints data table updates (rows 15-31) ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 25
$ws.Range("N15").Value = -82.758620689655
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = -12.5
$ws.Range("L16").Value = -9.090909090909
$ws.Range("M16").Value = -29.292929292929
$ws.Range("N16").Value = -85.59670781893
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 125
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 16
$ws.Range("I17").Value = 105
$ws.Range("J17").Value = 103
$ws.Range("K17").Value = 1.941747572815
$ws.Range("L17").Value = 10.526315789473
$ws.Range("M17").Value = 77.966101694915
$ws.Range("N17").Value = -60.674157303370
$ws.Range("C18").Value = 2
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 36
$ws.Range("K18").Value = 44
$ws.Range("L18").Value = -36.842105263157
$ws.Range("M18").Value = -43.75
$ws.Range("N18").Value = -94.321766561514
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 49
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 180
$ws.Range("J19").Value = 221
$ws.Range("K19").Value = -18.552036199095
$ws.Range("L19").Value = 5.882352941176
$ws.Range("M19").Value = 73.076923076923
$ws.Range("N19").Value = -50.138504155124
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 18.75
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = -12.5
$ws.Range("L20").Value = -22.222222222222
$ws.Range("M20").Value = 162.5
$ws.Range("N20").Value = -90.707964601769
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 27.586206896551
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = 6.140350877192
$ws.Range("I21").Value = 460
$ws.Range("J21").Value = 506
$ws.Range("K21").Value = -9.090909090909
$ws.Range("L21").Value = -4.761904761904
$ws.Range("M21").Value = 28.133704735376
$ws.Range("N21").Value = -81.459089076985
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("L22").Value = -68.75
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("C23").Value = 4
$ws.Range("F23").Value = 7
$ws.Range("I23").Value = 12
$ws.Range("K23").Value = 20
$ws.Range("L23").Value = 20
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -45.714285714285
$ws.Range("F24").Value = 93
$ws.Range("G24").Value = 110
$ws.Range("H24").Value = -15.454545454545
$ws.Range("I24").Value = 456
$ws.Range("J24").Value = 434
$ws.Range("K24").Value = 5.069124423963
$ws.Range("L24").Value = 16.326530612244
$ws.Range("M24").Value = 140
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -43.75
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 59
$ws.Range("H25").Value = -23.728813559322
$ws.Range("I25").Value = 207
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 3.5
$ws.Range("L25").Value = 30.188679245283
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 35.714285714285
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 28.947368421052
$ws.Range("I26").Value = 190
$ws.Range("J26").Value = 162
$ws.Range("K26").Value = 17.283950617283
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = 12.426035502958
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("I27").Value = 8
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 14.285714285714
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = -42.105263157894
$ws.Range("N29").Value = -94.366197183098
$ws.Range("N30").Value = -93.939393939393
$ws.Range("C14").Copy($ws.Range("F31"))
